$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player name, position, team) replacing the old rows 2-18
# and adding a new row 19.
$data = @(
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Moussa Diabate", "C", "Charlotte Hornets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Max Christie", "SG,SF", "Dallas Mavericks"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Jimmy Butler", "SF,PF", "Golden State Warriors"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Myles Turner", "C", "Indiana Pacers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
